$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (col E) and Correspond Handback DateTime (col H)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 08:45:22"
$wsZh.Range("H2").Value = "2016-03-19 08:46:02"
$wsZh.Range("E3").Value = "2016-03-19 08:45:22"
$wsZh.Range("H3").Value = "2016-03-19 08:46:02"

# de-de sheet: Correspond Handoff Datetime (col E) and Correspond Handback DateTime (col H)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 08:45:31"
$wsDe.Range("H2").Value = "2016-03-19 08:46:16"
$wsDe.Range("E3").Value = "2016-03-19 08:45:31"
$wsDe.Range("H3").Value = "2016-03-19 08:46:16"
